$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param(
        $ws,
        [string]$cellRef,
        [string]$text
    )
    $range = $ws.Range($cellRef)
    # Force a genuine text (shared-string) cell, byte-for-byte like typing
    # a value into a "General" formatted cell would, WITHOUT leaving a
    # leftover formula or a quotePrefix/NumberFormat style behind:
    #   1. write a formula that evaluates to the literal string
    #   2. copy it
    #   3. paste-special back over itself -> collapses to a plain static value
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial()
}

# Sheet tabs, in (1-based) tab order:
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha
# NB: sheet-name lookup (Worksheets.Item("...")) is case-insensitive, and
# "Vector_bf" / "Vector_BF" differ only by case, so Item("Vector_BF") would
# incorrectly resolve to the earlier "Vector_bf" tab. Index into
# $wb.Worksheets by (1-based) tab position instead, to stay unambiguous.

# --- Restricciones_del_lider ---
$ws = $wb.Worksheets.Item(2)
Set-TextValue $ws "A2" "1.0499999999999998 - x"
Set-TextValue $ws "B2" "-2.05"
Set-TextValue $ws "D2" "0.24"
Set-TextValue $ws "A3" "-1.05 + x"
Set-TextValue $ws "B3" "0.050000000000000044"
Set-TextValue $ws "D3" "0.72"

# --- Restricciones_del_follower ---
$ws = $wb.Worksheets.Item(3)
Set-TextValue $ws "A2" "-2.85 + y"
Set-TextValue $ws "B2" "1.85"
Set-TextValue $ws "D2" "0.47"
Set-TextValue $ws "E2" "1.3"
Set-TextValue $ws "F2" "1.4000000000000001"
Set-TextValue $ws "A3" "2.85 - y"
Set-TextValue $ws "B3" "-3.85"
Set-TextValue $ws "D3" "0.88"
Set-TextValue $ws "E3" "1.6"
Set-TextValue $ws "F3" "6.7"

# --- Punto_modificado ---
$ws = $wb.Worksheets.Item(4)
Set-TextValue $ws "A2" "1.05"
Set-TextValue $ws "B2" "2.85"

# --- Vector_bf ---
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws "A2" "-1.4248750000000001"

# --- Vector_BF ---
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws "A2" "-1.48"
Set-TextValue $ws "A3" "1.3"
